# frameData.xlsx edit
#
# The sheet contains two mirrored tables:
#   - a "live" table in rows 3:18 (columns B:G)
#   - a "base" table in rows 23:38 (columns B:G), which is the source of truth
#     for columns B,C,D,E; columns F and G in the base table are computed with
#     F = C-B+E-5 and G = D-B+E.
#
# The actual data edits are four cells in the base table:
#   B24: 12 -> 10
#   C26: 14 -> 13
#   C34: 26 -> 20
#   C37: 40 -> 30
#
# In addition, column F/G of the live table (rows 3:18), which previously held
# hard-coded numbers, are turned into formulas that simply mirror the base
# table's F/G values twenty rows below (F{r} = F{r+20}, G{r} = G{r+20}).
# Recalculation then updates every dependent cell (B4, C6, C8, C14, C17, etc.)
# automatically.
#
# Finally, the worksheet's saved cursor/selection moves from C29 to C27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Base table data edits (rows 23:38) ---
$ws.Range("B24").Value = 10
$ws.Range("C26").Value = 13
$ws.Range("C34").Value = 20
$ws.Range("C37").Value = 30

# --- Live table F/G columns become formulas mirroring the base table ---
$ws.Range("F3").Formula = "=F23"
$ws.Range("G3").Formula = "=G23"

$ws.Range("F4:G4").Formula = "=F24"
$ws.Range("F5:G5").Formula = "=F25"
$ws.Range("F6:G6").Formula = "=F26"
$ws.Range("F7:G7").Formula = "=F27"
$ws.Range("F8:G8").Formula = "=F28"
$ws.Range("F9:G9").Formula = "=F29"
$ws.Range("F10:G10").Formula = "=F30"
$ws.Range("F11:G11").Formula = "=F31"
$ws.Range("F12:G12").Formula = "=F32"
$ws.Range("F13:G13").Formula = "=F33"
$ws.Range("F14:G14").Formula = "=F34"
$ws.Range("F15:G15").Formula = "=F35"
$ws.Range("F16:G16").Formula = "=F36"
$ws.Range("F17:G17").Formula = "=F37"
$ws.Range("F18:G18").Formula = "=F38"

# --- Update the saved selection/active cell ---
$null = $ws.Range("C27").Select()
